$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" '63.207.19'
Set-TextValue $ws "E2" '  -5.00%  '
Set-TextValue $ws "D3" '3.302.89'
Set-TextValue $ws "E3" '  -5.91%  '
Set-TextValue $ws "D4" '1.00'
Set-TextValue $ws "E4" '  +0.08%  '
Set-TextValue $ws "D5" '549.07'
Set-TextValue $ws "E5" '  -2.26%  '
Set-TextValue $ws "D6" '169.85'
Set-TextValue $ws "E6" '  -8.52%  '
Set-TextValue $ws "D7" '0.606'
Set-TextValue $ws "E7" '  -4.22%  '
Set-TextValue $ws "E8" '  +0.16%  '
Set-TextValue $ws "D9" '3.288.01'
Set-TextValue $ws "E9" '  -6.21%  '
Set-TextValue $ws "D10" '0.610'
Set-TextValue $ws "E10" '  -4.80%  '
Set-TextValue $ws "D11" '0.148'
Set-TextValue $ws "E11" '  -5.12%  '
Set-TextValue $ws "D12" '52.93'
Set-TextValue $ws "E12" '  -4.58%  '
Set-TextValue $ws "D13" '0.0000262'
Set-TextValue $ws "E13" '  -6.39%  '
Set-TextValue $ws "D14" '8.85'
Set-TextValue $ws "E14" '  -5.67%  '
Set-TextValue $ws "D15" '3.837.67'
Set-TextValue $ws "E15" '  -5.84%  '
Set-TextValue $ws "E16" '  -4.20%  '
Set-TextValue $ws "D17" '3.319.47'
Set-TextValue $ws "E17" '  -5.58%  '
Set-TextValue $ws "D18" '17.60'
Set-TextValue $ws "E18" '  -6.20%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws "D19" '63.312.77'
Set-TextValue $ws "E19" '  -5.01%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws "D20" '11.56'
Set-TextValue $ws "E20" '  -4.30%  '
Set-TextValue $ws "D21" '0.961'
Set-TextValue $ws "E21" '  -4.24%  '
Set-TextValue $ws "D22" '402.23'
Set-TextValue $ws "D23" '4.03'
Set-TextValue $ws "E23" '  -1.70%  '
Set-TextValue $ws "D24" '4.24'
Set-TextValue $ws "E24" '  +1.86%  '
Set-TextValue $ws "D25" '82.22'
Set-TextValue $ws "E25" '  -4.12%  '
Set-TextValue $ws "D26" '13.05'
Set-TextValue $ws "E26" '  +5.59%  '
Set-TextValue $ws "D27" '10.52'
Set-TextValue $ws "E27" '  -4.55%  '
Set-TextValue $ws "D28" '2.70'
Set-TextValue $ws "E28" '  -7.40%  '
Set-TextValue $ws "D29" '8.51'
Set-TextValue $ws "E29" '  -7.07%  '
Set-TextValue $ws "D30" '28.89'
Set-TextValue $ws "E30" '  -5.10%  '
Set-TextValue $ws "D31" '6.43'
Set-TextValue $ws "E31" '  -4.03%  '
Set-TextValue $ws "D32" '573.26'
Set-TextValue $ws "E32" '  -8.10%  '
Set-TextValue $ws "D33" '11.22'
Set-TextValue $ws "E33" '  -5.35%  '
Set-TextValue $ws "D34" '0.105'
Set-TextValue $ws "E34" '  -6.12%  '
Set-TextValue $ws "D35" '57.25'
Set-TextValue $ws "E35" '  -4.54%  '
Set-TextValue $ws "D36" '0.999'
Set-TextValue $ws "E36" '  +0.05%  '
Set-TextValue $ws "E37" '  -2.69%  '
Set-TextValue $ws "D38" '34.77'
Set-TextValue $ws "E38" '  -9.07%  '
Set-TextValue $ws "D39" '3.38'
Set-TextValue $ws "E39" '  +0.77%  '
Set-TextValue $ws "D40" '0.0₃0729'
Set-TextValue $ws "E40" '  -11.42%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws "D41" '3.132.78'
Set-TextValue $ws "E41" '  -0.55%  '
$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws "D42" '0.363'
Set-TextValue $ws "E42" '  -6.50%  '
Set-TextValue $ws "E43" '  +0.03%  '
Set-TextValue $ws "D44" '2.78'
Set-TextValue $ws "E44" '  -3.32%  '
Set-TextValue $ws "D45" '3.15'
Set-TextValue $ws "E45" '  -4.01%  '
Set-TextValue $ws "D46" '2.42'
Set-TextValue $ws "E46" '  -8.76%  '
Set-TextValue $ws "D47" '0.0399'
Set-TextValue $ws "E47" '  -4.63%  '
Set-TextValue $ws "D48" '2.58'
Set-TextValue $ws "E48" '  -4.93%  '
Set-TextValue $ws "E49" '  -4.98%  '
Set-TextValue $ws "D50" '132.63'
Set-TextValue $ws "E50" '  -5.16%  '
Set-TextValue $ws "D51" '7.95'
Set-TextValue $ws "E51" '  -6.75%  '
